$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.029136180877686
$ws.Range("B1").Value = 3.27988600730896
$ws.Range("C1").Value = 3.677816867828369
$ws.Range("D1").Value = 2.004348039627075
$ws.Range("E1").Value = 1.175016760826111
